# Weekly refresh of the "Feria Lagunitas de Puerto Montt - Mango" series:
# a new daily record is inserted at the top of the data block (row 232),
# pushing the existing rows 232:274 down to 233:275.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row, shifting rows 232:274 -> 233:275 and extending
# the sheet's used range / dimension to A1:T275 automatically.
$ws.Rows(232).Insert()

# Populate the new row 232 with the latest price observation.
$ws.Range("A232").Value = 4
$ws.Range("B232").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C232").Value = "Los Lagos"
$ws.Range("D232").Value = 44889
$ws.Range("E232").Value = 10
$ws.Range("F232").Value = "Fruta"
$ws.Range("G232").Value = 100108
$ws.Range("H232").Value = "Tropicales y subtropicales"
$ws.Range("I232").Value = 100108002
$ws.Range("J232").Value = "Mango"
$ws.Range("K232").Value = "Sin especificar"
$ws.Range("L232").Value = "Primera"
$ws.Range("M232").Value = 80
$ws.Range("N232").Value = 9000
$ws.Range("O232").Value = 10000
$ws.Range("P232").Value = 9500
$ws.Range("Q232").Value = "`$/bandeja 4 kilos"
$ws.Range("R232").Value = "Brasil"
$ws.Range("S232").Value = 2375
$ws.Range("T232").Value = 4
